# Update the threshold values in Sheet1 to reflect the re-uploaded data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (alpha_distance_range): Min 5.5 -> 5.6, Max 10 -> 9
$ws.Range("B2").Value = 5.6
$ws.Range("C2").Value = 9

# Row 3 (beta_distance_range): Min 6.5 -> 6.6, Max 8.5 -> 8.3
$ws.Range("B3").Value = 6.6
$ws.Range("C3").Value = 8.3

# Row 4 (ratio_threshold_range): Min 0.7 -> 0.75
$ws.Range("B4").Value = 0.75

$wb.Save()
